$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.002679
$ws.Range("H2").Value = 0.116782
$ws.Range("I2").Value = 0.419822
$ws.Range("J2").Value = 1.4807
$ws.Range("K2").Value = -0.700134
$ws.Range("L2").Value = -0.125914
$ws.Range("M2").Value = -1.425507
$ws.Range("N2").Value = -0.784363
$ws.Range("O2").Value = -0.123064
$ws.Range("P2").Value = -0.057872
$ws.Range("Q2").Value = 1.367716
$ws.Range("R2").Value = -0.170843
$ws.Range("G3").Value = 0.011203
$ws.Range("H3").Value = 0.224431
$ws.Range("I3").Value = 0.346664
$ws.Range("J3").Value = -1.421002
$ws.Range("K3").Value = -0.749111
$ws.Range("L3").Value = -0.084967
$ws.Range("M3").Value = 1.473656
$ws.Range("N3").Value = -0.6730699999999999
$ws.Range("O3").Value = -0.091923
$ws.Range("P3").Value = -0.063856
$ws.Range("Q3").Value = 1.19775
$ws.Range("R3").Value = -0.169773
$ws.Range("G4").Value = -0.129309
$ws.Range("H4").Value = 0.00706
$ws.Range("I4").Value = 0.370343
$ws.Range("J4").Value = 1.590831
$ws.Range("K4").Value = -0.043715
$ws.Range("L4").Value = -0.095932
$ws.Range("M4").Value = -0.749094
$ws.Range("N4").Value = -1.115596
$ws.Range("O4").Value = -0.136415
$ws.Range("P4").Value = -0.712429
$ws.Range("Q4").Value = 1.152252
$ws.Range("R4").Value = -0.137996
$ws.Range("G5").Value = -0.14432
$ws.Range("H5").Value = 0.09843
$ws.Range("I5").Value = 0.341551
$ws.Range("J5").Value = 1.549684
$ws.Range("K5").Value = -0.254352
$ws.Range("L5").Value = -0.080599
$ws.Range("M5").Value = -0.968772
$ws.Range("N5").Value = -0.901166
$ws.Range("O5").Value = -0.113244
$ws.Range("P5").Value = -0.436591
$ws.Range("Q5").Value = 1.057088
$ws.Range("R5").Value = -0.147709
$ws.Range("G6").Value = -0.209128
$ws.Range("H6").Value = 0.019171
$ws.Range("I6").Value = 0.327533
$ws.Range("J6").Value = 1.523963
$ws.Range("K6").Value = -0.192741
$ws.Range("L6").Value = -0.06741800000000001
$ws.Range("M6").Value = -0.537409
$ws.Range("N6").Value = 0.978013
$ws.Range("O6").Value = -0.128213
$ws.Range("P6").Value = -0.777426
$ws.Range("Q6").Value = -0.804442
$ws.Range("R6").Value = -0.131902
$ws.Range("G7").Value = 0.00177
$ws.Range("H7").Value = 0.009312000000000001
$ws.Range("I7").Value = 0.353032
$ws.Range("J7").Value = -0.662648
$ws.Range("K7").Value = -1.117848
$ws.Range("L7").Value = -0.115516
$ws.Range("M7").Value = -0.637657
$ws.Range("N7").Value = 1.127546
$ws.Range("O7").Value = -0.119433
$ws.Range("P7").Value = 1.298535
$ws.Range("Q7").Value = -0.019009
$ws.Range("R7").Value = -0.118083
$ws.Range("G8").Value = -0.004847
$ws.Range("H8").Value = 0.08786099999999999
$ws.Range("I8").Value = 0.333904
$ws.Range("J8").Value = 1.16119
$ws.Range("K8").Value = -0.533231
$ws.Range("L8").Value = -0.099868
$ws.Range("M8").Value = -1.084401
$ws.Range("N8").Value = -0.647656
$ws.Range("O8").Value = -0.099506
$ws.Range("P8").Value = -0.071941
$ws.Range("Q8").Value = 1.093026
$ws.Range("R8").Value = -0.13453
$ws.Range("G9").Value = -0.09139600000000001
$ws.Range("H9").Value = 0.004762
$ws.Range("I9").Value = 0.312238
$ws.Range("J9").Value = 1.227146
$ws.Range("K9").Value = -0.047013
$ws.Range("L9").Value = -0.0829
$ws.Range("M9").Value = -0.5911920000000001
$ws.Range("N9").Value = -0.884832
$ws.Range("O9").Value = -0.114385
$ws.Range("P9").Value = -0.544558
$ws.Range("Q9").Value = 0.927084
$ws.Range("R9").Value = -0.114953
$ws.Range("G10").Value = 0.010979
$ws.Range("H10").Value = -0.007273
$ws.Range("I10").Value = 0.277261
$ws.Range("J10").Value = 1.026126
$ws.Range("K10").Value = -0.237963
$ws.Range("L10").Value = -0.095317
$ws.Range("M10").Value = -0.734105
$ws.Range("N10").Value = -0.751586
$ws.Range("O10").Value = -0.091922
$ws.Range("P10").Value = -0.303
$ws.Range("Q10").Value = 0.996822
$ws.Range("R10").Value = -0.090022
